$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting B:H (formerly A:G) to the right
$ws.Columns.Item(1).Insert()

# Header for the new column, matching the style used by the other header cells
$ws.Range("A1").Value = "CLAVE DE LA ENTIDAD"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Fill column A rows 2-33 with sequential numeric keys 1-32
for ($i = 2; $i -le 33; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}
